# Fruta / hortaliza, semanal
#
# This sheet keeps a rolling weekly history: the newest reading always
# lands in row 2, and everything that was already there shifts down one
# row. Here a brand-new "Primera" quality record (Provincia de Quillota)
# is added as the new row 2; the former rows 2 and 3 move down to become
# rows 3 and 4 unchanged.
#
# We write the rows from the bottom up using literal values (no Range-to-
# Range .Value copy, no Rows.Insert()) so no incidental formatting/styles
# get pulled in along the way - only column D keeps its original
# "YYYY-MM-DD HH:MM:SS" date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- what used to be row 3 (La Ligua, Primera, 80 vol.)
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44466
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104004
$ws.Range("J4").Value = "Níspero"
$ws.Range("K4").Value = "Californiana(o)"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 11000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 11000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("R4").Value = "La Ligua"
$ws.Range("S4").Value = 2200
$ws.Range("T4").Value = 5

# Row 3 <- what used to be row 2 (La Ligua, Segunda, 20 vol.)
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44166
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104004
$ws.Range("J3").Value = "Níspero"
$ws.Range("K3").Value = "Californiana(o)"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "La Ligua"
$ws.Range("S3").Value = 667
$ws.Range("T3").Value = 18

# Row 2 <- this week's new record (Provincia de Quillota, Primera, 35 vol.)
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44483
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104004
$ws.Range("J2").Value = "Níspero"
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 10000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 2000
$ws.Range("T2").Value = 5
